# Updated cryptos list values (prices / 1h volume %) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.536.02'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.844.57'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''262.58'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = '''0.5328'
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('D8').Value = '''0.3125'
$ws.Range('E8').Value = '  -4.66%  '
$ws.Range('D9').Value = '''0.06899'
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').Value = '''18.67'
$ws.Range('E10').Value = '  -0.98%  '
$ws.Range('D11').Value = '''0.7650'
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').Value = '''0.07836'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '1.840.41'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '''89.68'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = '''5.046'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '''1.000'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '''14.06'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.000007957'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '''1.001'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').Value = '26.559.84'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').Value = '2.081.09'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').Value = '''4.631'
$ws.Range('D23').Value = '''6.023'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '''9.323'
$ws.Range('E24').Value = '  -2.24%  '
$ws.Range('D25').Value = '''141.50'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').Value = '''2.191'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').Value = '''1.692'
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('D28').Value = '''17.05'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '''111.30'
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('D30').Value = '''4.285'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').Value = '''0.08801'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('D32').Value = '''4.093'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').Value = '''0.04840'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''0.7367'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.941'
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('D36').Value = '''1.137'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').Value = '''2.335'
$ws.Range('E38').Value = '  +5.79%  '
$ws.Range('D39').Value = '''0.01728'
$ws.Range('E39').Value = '  -3.38%  '
$ws.Range('D40').Value = '''0.4822'
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('D41').Value = '''0.9053'
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').Value = '''108.49'
$ws.Range('E42').Value = '  -3.96%  '
$ws.Range('D43').Value = '''5.904'
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '''7.660'
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('D46').Value = '''0.4166'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('D47').Value = '''9.011'
$ws.Range('E47').Value = '  -1.01%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').Value = '''35.00'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').Value = '''0.05812'
$ws.Range('E50').Value = '  -2.00%  '
$ws.Range('D51').Value = '''0.8973'
$ws.Range('E51').Value = '  +1.02%  '
